$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "kavanan"
$ws.Range("A3").Value = "sharat"
$ws.Range("A1").Value = "hi"

$ws.Range("A2").Select()
